$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105, pushing the existing rows 105-140 down to 106-141
$ws.Rows(105).Insert()

# Populate the newly inserted row 105 with the new weekly record
$ws.Cells.Item(105, 1).Value  = 7
$ws.Cells.Item(105, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(105, 3).Value  = "Ñuble"
$ws.Cells.Item(105, 4).Value  = 45229
$ws.Cells.Item(105, 5).Value  = 16
$ws.Cells.Item(105, 6).Value  = 100112013
$ws.Cells.Item(105, 7).Value  = "Alcachofa"
$ws.Cells.Item(105, 8).Value  = "Española"
$ws.Cells.Item(105, 9).Value  = "Primera"
$ws.Cells.Item(105, 10).Value = 80
$ws.Cells.Item(105, 11).Value = 12000
$ws.Cells.Item(105, 12).Value = 12000
$ws.Cells.Item(105, 13).Value = 12000
$ws.Cells.Item(105, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(105, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(105, 16).Value = 400
$ws.Cells.Item(105, 17).Value = 30
$ws.Cells.Item(105, 18).Value = "Hortaliza"
